# Append the new resale-number row (2023-06-10 21:18:14) to the
# CityResaleNum sheet, mirroring the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

# Text columns: prefix with a leading apostrophe so Excel stores them as
# literal text instead of auto-converting to a date/time/number.
$ws.Cells.Item($row, 1).Value = "'2023-06-10"
$ws.Cells.Item($row, 2).Value = "'21:18:14"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "'23"

# Numeric columns (Beijing .. Wuhan).
$ws.Cells.Item($row, 5).Value  = 121201
$ws.Cells.Item($row, 6).Value  = 134643
$ws.Cells.Item($row, 7).Value  = 161000
$ws.Cells.Item($row, 8).Value  = 132032
$ws.Cells.Item($row, 9).Value  = 176211
$ws.Cells.Item($row, 10).Value = 114167
$ws.Cells.Item($row, 11).Value = 202009
$ws.Cells.Item($row, 12).Value = 222446
$ws.Cells.Item($row, 13).Value = 173715
$ws.Cells.Item($row, 14).Value = 98983
$ws.Cells.Item($row, 15).Value = 38774
$ws.Cells.Item($row, 16).Value = 34281
$ws.Cells.Item($row, 17).Value = 51203
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 37137
$ws.Cells.Item($row, 20).Value = -1
